$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, centered alignment) from the last
# existing header cell (AB1) onto the three new header cells, then set their
# text values.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 86  # AC - Wins
    $ws.Cells.Item($r, 30).Value = 76  # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0   # AE - Ties
}
